# statystyka błędów Carmen + nowe struktury przesyłane do GUI - zbyt duże MCU LOAD
# Fills rows 42-48 of the "Translation" sheet with new text-id / typography /
# alignment / direction / translation rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# TEXT ID, TYPOGRAPHY NAME, ALIGNMENT, DIRECTION, GB
$rows = @(
    @("SingleUseId61", "Typography_00", "Center", "LTR", "Temperature [C]:", $false),
    @("SingleUseId62", "Typography_00", "Center", "LTR", "Valid / E_CRC:", $false),
    @("SingleUseId63", "Typography_00", "Left",   "LTR", "<value>", $false),
    @("SingleUseId64", "Typography_00", "Left",   "LTR", "0,00", $false),
    @("SingleUseId65", "Typography_00", "Center", "LTR", "<value>/<value>", $false),
    @("SingleUseId67", "Typography_00", "Left",   "LTR", "10000000", $true),
    @("SingleUseId68", "Typography_00", "Left",   "LTR", "1000000", $true)
)

$startRow = 42
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 2).Value = $data[0]
    $ws.Cells.Item($r, 3).Value = $data[1]
    $ws.Cells.Item($r, 4).Value = $data[2]
    $ws.Cells.Item($r, 5).Value = $data[3]

    $fCell = $ws.Cells.Item($r, 6)
    $needsTextFormat = $data[5]
    if ($needsTextFormat) {
        # Purely-numeric-looking translation text (e.g. "10000000") must stay
        # text, not get auto-converted to a number by Excel.
        $fCell.NumberFormat = "@"
        $fCell.Value = $data[4]
        $fCell.Style = "Normal"
    } else {
        $fCell.Value = $data[4]
    }
}
